$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# (e.g. "531.72", "1.00", "0.0930") are not coerced into floating-point
# numbers and lose their original formatting/precision.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '57.033.74'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.348.08'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '531.72'
$ws.Range('E5').Value = '  +2.05%  '
$ws.Range('D6').Value = '132.70'
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.535'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').Value = '2.347.09'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').Value = '2.741.77'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '23.48'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('D16').Value = '57.078.41'
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = '2.337.05'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = '337.55'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '61.75'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '8.73'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '0.994'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').Value = '174.23'
$ws.Range('E29').Value = '  +3.49%  '
$ws.Range('D30').Value = '1.72'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('E31').Value = '  -2.79%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = '18.53'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '0.994'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -4.35%  '
$ws.Range('D37').Value = '0.914'
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = '39.27'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').Value = '1.57'
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '5.78'
$ws.Range('E41').Value = '  +7.72%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '149.33'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').Value = '282.38'
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').Value = '0.0930'
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('D47').Value = '0.0501'
$ws.Range('E47').Value = '  -1.86%  '
$ws.Range('D48').Value = '18.89'
$ws.Range('E48').Value = '  +3.06%  '
$ws.Range('D49').Value = '0.559'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -1.69%  '
$ws.Range('E51').Value = '  -2.24%  '

# Restore the default (unstyled) cell style on column D now that the
# text values are committed, so no stray style index is left behind.
$ws.Range('D2:D51').Style = 'Normal'

